# Update COVID recession files for July STEO
#
# The "Data" worksheet holds EIA Short-Term Energy Outlook (STEO) GDP
# figures. This edit refreshes the data from the "May STEO" release to
# the "July STEO" release, updating the 2020 and 2021 GDP projections.
# All other cells on this and other sheets are formulas that will
# recalculate automatically.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Data")

# Label: "May STEO" -> "July STEO"
$ws.Range("A3").Value = "July STEO"

# Updated STEO projections for 2020 (C3) and 2021 (D3); 2019 actual (B3)
# is unchanged.
$ws.Range("C3").Value = 17517
$ws.Range("D3").Value = 18418

$excel.CalculateFullRebuild()
